$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 82: date serial 45884 (2025-08-15) in A82, value 105.26 in B82
$ws.Cells.Item(82, 1).Value = 45884
$ws.Cells.Item(82, 2).Value = 105.26

# Match formatting of the row above (A81 has a date-style format),
# copying the style index (s="2") rather than creating a new style entry.
$ws.Cells.Item(81, 1).Copy()
$ws.Cells.Item(82, 1).PasteSpecial(-4122)
